$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4488
$ws1.Range("F3").Value = 2480
$ws1.Range("F5").Value = 30
$ws1.Range("F10").Value = 167
$ws1.Range("F12").Value = 1676
$ws1.Range("F13").Value = 303
$ws1.Range("F14").Value = 3646
$ws1.Range("F15").Value = 13

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4488
$ws4.Range("F3").Value = 2480
$ws4.Range("F5").Value = 30
$ws4.Range("F12").Value = 167
$ws4.Range("F16").Value = 1676
$ws4.Range("F17").Value = 303
$ws4.Range("F18").Value = 3646
$ws4.Range("F19").Value = 13
